# Fixes #1949 - data importer parsing unit should be case insensitive.
#
# Test5.xlsx is the fixture used by the data importer integration tests.
# The unit text in B2 is changed from the lower-case "mg/l" to the
# mixed-case "mg/mL" so the round trip exercises case-insensitive unit
# parsing. The rest of the edits mirror the small formatting touch-ups
# left behind by the Excel session that produced the fixture (selection
# location and column B being widened to fit the new text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2 ("H" row) holds the unit string - switch it to the mixed-case form.
$ws.Range("B2").Value = "mg/mL"

# Column B is widened so "mg/mL" / "Concentration" are fully visible.
$ws.Columns.Item(2).ColumnWidth = 13.428571428571429

# Leave the selection where the editing session left it.
[void]$ws.Range("D9").Select()

# The built-in "Normal" cell style was saved under its German display
# name ("Standard"); rename it back to the canonical "Normal" name.
$wb.Styles.Item("Standard").Delete()
$wb.Styles.Add("Normal")
